$d = $word.ActiveDocument

$replacements = @(
    @{old="87÷4="; new="98÷6="},
    @{old="55÷5="; new="29÷3="},
    @{old="83÷4="; new="32÷8="},
    @{old="85÷3="; new="41÷6="},
    @{old="70÷8="; new="30÷4="},
    @{old="89÷8="; new="59÷8="},
    @{old="70÷6="; new="95÷6="},
    @{old="74÷6="; new="26÷5="},
    @{old="61÷8="; new="12÷7="},
    @{old="41÷8="; new="16÷5="},
    @{old="17÷7="; new="14÷8="},
    @{old="31÷6="; new="69÷2="},
    @{old="58÷6="; new="13÷2="},
    @{old="40÷2="; new="13÷7="},
    @{old="57÷6="; new="34÷2="},
    @{old="68÷9="; new="48÷5="},
    @{old="54÷6="; new="80÷9="},
    @{old="64÷2="; new="72÷3="},
    @{old="98÷2="; new="39÷3="},
    @{old="75÷4="; new="31÷6="},
    @{old="15÷7="; new="37÷4="},
    @{old="46÷2="; new="10÷3="},
    @{old="54÷7="; new="19÷9="},
    @{old="12÷8="; new="41÷9="},
    @{old="94÷6="; new="57÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
